$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A52").Value = "frontend form articulos el codigo solo en mayusculas"
$ws.Range("B52").Value = "no comenzado"

$ws.Range("A48").Select()
